$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "4913" course-equivalence column (AG) — it held
# ECE4913 (row1) and EE4913 (row2) which are being retired from the list.
$ws.Range("AG1:AG2").ClearContents()

# Add a brand-new column (AI) for the master list GUI feature with
# placeholder values a / b / c / y on rows 1-4.
$ws.Range("AI1").Value = "a"
$ws.Range("AI2").Value = "b"
$ws.Range("AI3").Value = "c"
$ws.Range("AI4").Value = "y"
